$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update header timestamp text (A1)
$ws.Range("A1").Value = "Datos actualizados a 19 de Mayo de 2020 a las 16:35"

# Update Estados Unidos (row 4) totals
$ws.Range("B4").Value = 1552315
$ws.Range("C4").Value = 2021
$ws.Range("D4").Value = 358923
$ws.Range("E4").Value = 1101320

# Reorder Belice / Nueva Caledonia (rows 196-197) along with their stats
$ws.Range("A196").Value = "Belice"
$ws.Range("D196").Value = 16
$ws.Range("H196").Value = 2

$ws.Range("A197").Value = "Nueva Caledonia"
$ws.Range("D197").Value = 18
$ws.Range("H197").Value = 0

# Reorder Groenlandia / Montserrat / Seychelles (rows 209-211) along with their stats
$ws.Range("A209").Value = "Groenlandia"

$ws.Range("A210").Value = "Montserrat"
$ws.Range("D210").Value = 10
$ws.Range("H210").Value = 1

$ws.Range("A211").Value = "Seychelles"
$ws.Range("D211").Value = 11
$ws.Range("H211").Value = 0

# Reorder San Bartolome / Bonaire, San Eustaquio y Saba (rows 215-216)
$ws.Range("A215").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("A216").Value = "San Bartolome"
